$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new trade row (row 7) mirroring the existing data rows.
$ws.Range("A7").Value = 10092.81
$ws.Range("B7").Value = 10156.799999999999
$ws.Range("C7").Value = 307.87
$ws.Range("D7").Value = 305.93
$ws.Range("E7").Value = $false
$ws.Range("F7").Value = -0.63
$ws.Range("G7").Value = 42612.675358796296
$ws.Range("G7").NumberFormat = "m/d/yy h:mm"
$ws.Range("H7").Value = $false
